$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns AD/AE/AF -> Wins / Losses / Ties
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy style of existing header cell (AC1) to the new header cells so they
# match the bold/centered/bordered header formatting.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null

# Data rows 2-50: team record values are constant per the source data.
for ($row = 2; $row -le 50; $row++) {
    $ws.Cells.Item($row, 30).Value = 87
    $ws.Cells.Item($row, 31).Value = 75
    $ws.Cells.Item($row, 32).Value = 0
}
